# Update TPM-derived LR-pair statistics (Mdk-Itga4) with refreshed expression values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.01253823027
$ws.Range("R2").Value = 0.11284407243
$ws.Range("S2").Value = 0.0003139206594611955
$ws.Range("T2").Value = 0.0003139206594611955

# Row 3
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.07156975043666665
$ws.Range("R3").Value = 0.6441277539299999
$ws.Range("S3").Value = 0.001791897482398971
$ws.Range("T3").Value = 0.001791897482398972

# Row 4
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 1.262537938806667
$ws.Range("R4").Value = 11.36284144926
$ws.Range("S4").Value = 0.0316102618798823
$ws.Range("T4").Value = 0.0316102618798823

# Row 5
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.2803734452180001
$ws.Range("R5").Value = 2.523361006962
$ws.Range("S5").Value = 0.007019732045345658
$ws.Range("T5").Value = 0.007019732045345658

# Row 6
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.04006948826101438
$ws.Range("T6").Value = 0.04006948826101438

# Row 7
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 28.23222288943156
$ws.Range("R7").Value = 254.0900060048841
$ws.Range("S7").Value = 0.7068523895841494
$ws.Range("T7").Value = 0.7068523895841493

# Row 8
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.07896519341
$ws.Range("R8").Value = 0.7106867406899999
$ws.Range("S8").Value = 0.001977057770988604
$ws.Range("T8").Value = 0.001977057770988605

# Row 9
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.4507429727988888
$ws.Range("R9").Value = 4.056686755189999
$ws.Range("S9").Value = 0.01128528733493478
$ws.Range("T9").Value = 0.01128528733493478

# Row 10
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 7.951405451842223
$ws.Range("R10").Value = 71.56264906658001
$ws.Range("S10").Value = 0.1990799649818248
$ws.Range("T10").Value = 0.1990799649818248
